# Applies the commit: "Added code to remove unnecessary columns and created
# final dataframe to be used downstream"
#
# 1. Hides the 30 rows in "fields" that are marked (red-highlighted) for
#    removal, and turns on a cell-color AutoFilter on column A so the sheet
#    keeps its "filtered" state.
# 2. Adds a new worksheet "fields to keep" after "fields" containing the 42
#    field names that were kept, plus a helper column D that builds up a
#    Python list literal (e.g. "['basements',", "'zoning']") from column A.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. "fields" sheet: turn on a filter-by-color on column A and hide the
#    rows that correspond to the red-highlighted (removed) fields.
# ---------------------------------------------------------------------
$sampleColoredCell = $ws1.Range("A2")
$ws1.Range("A1:D73").AutoFilter(1, $sampleColoredCell, 8)

$hiddenRowNums = @(2,4,5,7,11,12,22,24,25,28,29,30,31,32,33,40,41,42,43,44,46,47,50,51,53,54,55,59,66,68)
for ($r = 2; $r -le 73; $r++) {
    $ws1.Rows.Item($r).Hidden = $false
}
foreach ($r in $hiddenRowNums) {
    $ws1.Rows.Item($r).Hidden = $true
}

$app = $wb.Application
$app.ActiveWindow.ScrollRow = 16
$ws1.Range("A69").Select()

# ---------------------------------------------------------------------
# 2. Add the "fields to keep" sheet right after "fields".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "fields to keep"

$fields = @("basements","building_code_description","category_code_description","census_tract","central_air","depth","exempt_building","exempt_land","exterior_condition","fireplaces","frontage","fuel","garage_spaces","garage_type","geographic_ward","interior_condition","location","market_value","market_value_date","number_of_bathrooms","number_of_bedrooms","number_of_rooms","number_stories","quality_grade","sale_date","sale_price","shape","street_designation","street_direction","street_name","taxable_building","taxable_land","topography","total_area","total_livable_area","type_heater","unit","view","year_built","year_built_estimate","zip_code","zoning")

for ($i = 1; $i -le $fields.Count; $i++) {
    $newSheet.Range("A" + $i).Value = $fields[$i - 1]
    if ($i -eq 1) {
        $newSheet.Range("D" + $i).Formula = '="[''"&A' + $i + '&"'',"'
    } elseif ($i -eq $fields.Count) {
        $newSheet.Range("D" + $i).Formula = '="''"&A' + $i + '&"'']"'
    } else {
        $newSheet.Range("D" + $i).Formula = '="''"&A' + $i + '&"'',"'
    }
}

# B1 carries a leftover "cleared formatting" style in the source workbook.
$newSheet.Range("B1").Value = 0
$newSheet.Range("B1").NumberFormat = "General"
$newSheet.Range("B1").ClearContents()

$newSheet.Range("D1:D42").Select()
$app.ActiveWindow.ScrollRow = 23

$ws1.Activate()
